$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 747, shifting existing rows 747:792 down to 748:793
$ws.Rows.Item(747).Insert()

# Populate the newly inserted row 747 with the new data record
$ws.Cells.Item(747, 1).Value = 5
$ws.Cells.Item(747, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(747, 3).Value = "Maule"
$ws.Cells.Item(747, 4).Value = 45147
$ws.Cells.Item(747, 5).Value = 7
$ws.Cells.Item(747, 6).Value = "Fruta"
$ws.Cells.Item(747, 7).Value = 100109
$ws.Cells.Item(747, 8).Value = "Uva"
$ws.Cells.Item(747, 9).Value = 100109001
$ws.Cells.Item(747, 10).Value = "Uva"
$ws.Cells.Item(747, 11).Value = "Crimpson Seedless"
$ws.Cells.Item(747, 12).Value = "Segunda"
$ws.Cells.Item(747, 13).Value = 250
$ws.Cells.Item(747, 14).Value = 12000
$ws.Cells.Item(747, 15).Value = 12000
$ws.Cells.Item(747, 16).Value = 12000
$ws.Cells.Item(747, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(747, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(747, 19).Value = 1500
$ws.Cells.Item(747, 20).Value = 8

Write-Host "Done"
